$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17 (hunk 0)
$ws.Range("H17").Value = 998.7547
$ws.Range("I17").Value = 572.6923
$ws.Range("J17").Value = 1409.037
$ws.Range("K17").Value = 1718.0769
$ws.Range("L17").Value = 4227.111
$ws.Range("M17").Value = -1550.0769
$ws.Range("N17").Value = -4563.111
# Row 34 (hunk 1)
$ws.Range("H34").Value = 7661.7144
$ws.Range("I34").Value = 605.3333
$ws.Range("J34").Value = 50000
$ws.Range("K34").Value = 605.3333
$ws.Range("L34").Value = 50000
$ws.Range("M34").Value = -402.3333
$ws.Range("N34").Value = -50406
# Row 36 (hunk 2)
$ws.Range("H36").Value = 7661.7144
$ws.Range("I36").Value = 605.3333
$ws.Range("J36").Value = 50000
$ws.Range("K36").Value = 605.3333
$ws.Range("L36").Value = 50000
$ws.Range("M36").Value = 109.6667
$ws.Range("N36").Value = -51430
# Row 107 (hunk 3)
$ws.Range("H107").Value = 1042.5385
$ws.Range("I107").Value = 1109.1305
$ws.Range("J107").Value = 532
$ws.Range("K107").Value = 1109.1305
$ws.Range("L107").Value = 532
$ws.Range("M107").Value = 810.8695
$ws.Range("N107").Value = -4372
# Row 115 (hunk 4)
$ws.Range("H115").Value = 3101.875
$ws.Range("I115").Value = 2002.5
$ws.Range("J115").Value = 6400
$ws.Range("K115").Value = 6007.5
$ws.Range("L115").Value = 19200
$ws.Range("M115").Value = -4440.5
$ws.Range("N115").Value = -22334
# Row 137 (hunk 5)
$ws.Range("H137").Value = 1246.4529
$ws.Range("I137").Value = 1137.4117
$ws.Range("J137").Value = 1441.579
$ws.Range("K137").Value = 3412.2351
$ws.Range("L137").Value = 4324.737
$ws.Range("M137").Value = -862.2351000000003
$ws.Range("N137").Value = -9424.737000000001
# Row 138 (hunk 6)
$ws.Range("H138").Value = 2087.6265
$ws.Range("I138").Value = 1431.0834
$ws.Range("J138").Value = 2820.5117
$ws.Range("K138").Value = 4293.2502
$ws.Range("L138").Value = 8461.535100000001
$ws.Range("M138").Value = 846.7497999999996
$ws.Range("N138").Value = -18741.5351
# Row 141 (hunk 7)
$ws.Range("H141").Value = 2485.9343
$ws.Range("I141").Value = 749.3913
$ws.Range("J141").Value = 7811.3335
$ws.Range("K141").Value = 2248.1739
$ws.Range("L141").Value = 23434.0005
$ws.Range("M141").Value = 2931.8261
$ws.Range("N141").Value = -33794.00049999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61 (hunk 8)
$ws.Range("H61").Value = 2117.5942
$ws.Range("I61").Value = 1779.7084
$ws.Range("J61").Value = 2889.9048
$ws.Range("K61").Value = 1779.7084
$ws.Range("L61").Value = 2889.9048
$ws.Range("M61").Value = -1567.7084
$ws.Range("N61").Value = -3313.9048
# Row 74 (hunk 9)
$ws.Range("I74").Value = 671.825
$ws.Range("J74").Value = 1341.6666
$ws.Range("K74").Value = 671.825
$ws.Range("L74").Value = 1341.6666
$ws.Range("M74").Value = 202.175
$ws.Range("N74").Value = -3089.6666
# Row 77 (hunk 10)
$ws.Range("I77").Value = 671.825
$ws.Range("J77").Value = 1341.6666
$ws.Range("K77").Value = 3359.125
$ws.Range("L77").Value = 6708.333000000001
$ws.Range("M77").Value = 1008.875
$ws.Range("N77").Value = -15444.333
# Row 97 (hunk 11)
$ws.Range("H97").Value = 826.5185
$ws.Range("I97").Value = 735
$ws.Range("J97").Value = 1352.75
$ws.Range("K97").Value = 735
$ws.Range("L97").Value = 1352.75
$ws.Range("M97").Value = -239
$ws.Range("N97").Value = -2344.75
# Row 132 (hunk 12)
$ws.Range("H132").Value = 3784
$ws.Range("I132").Value = 3467.647
$ws.Range("J132").Value = 4416.706
$ws.Range("K132").Value = 10402.941
$ws.Range("L132").Value = 13250.118
$ws.Range("M132").Value = -7872.940999999999
$ws.Range("N132").Value = -18310.118
# Row 136 (hunk 13)
$ws.Range("H136").Value = 2117.5942
$ws.Range("I136").Value = 1779.7084
$ws.Range("J136").Value = 2889.9048
$ws.Range("K136").Value = 5339.1252
$ws.Range("L136").Value = 8669.714399999999
$ws.Range("M136").Value = -2789.1252
$ws.Range("N136").Value = -13769.7144

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 57 (hunk 14)
$ws.Range("H57").Value = 13732.381
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 13732.381
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 13732.381
$ws.Range("N57").Value = -15172.381
# Row 105 (hunk 15)
$ws.Range("H105").Value = 41669330
$ws.Range("I105").Value = 41669330
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 41669330
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -41667583
# Row 107 (hunk 16)
$ws.Range("H107").Value = 78264.46000000001
$ws.Range("I107").Value = 112102.78
$ws.Range("J107").Value = 2128.25
$ws.Range("K107").Value = 112102.78
$ws.Range("L107").Value = 2128.25
$ws.Range("M107").Value = -110182.78
$ws.Range("N107").Value = -5968.25
# Row 132 (hunk 17)
$ws.Range("H132").Value = 38108.9
$ws.Range("I132").Value = 25354.5
$ws.Range("J132").Value = 41297.5
$ws.Range("K132").Value = 25354.5
$ws.Range("L132").Value = 41297.5
$ws.Range("M132").Value = -20294.5
$ws.Range("N132").Value = -51417.5
# Row 133 (hunk 18)
$ws.Range("H133").Value = 55780
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 55780
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 55780
$ws.Range("N133").Value = -65900
# Row 134 (hunk 19)
$ws.Range("H134").Value = 2736.0715
$ws.Range("I134").Value = 2516.04
$ws.Range("J134").Value = 3059.647
$ws.Range("K134").Value = 7548.12
$ws.Range("L134").Value = 9178.940999999999
$ws.Range("M134").Value = -5013.12
$ws.Range("N134").Value = -14248.941
# Row 135 (hunk 20)
$ws.Range("H135").Value = 49390
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 49390
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 49390
$ws.Range("N135").Value = -59530
# Row 136 (hunk 21)
$ws.Range("H136").Value = 13732.381
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 13732.381
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 13732.381
$ws.Range("N136").Value = -23932.381

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (hunk 22)
$ws.Range("H31").Value = 4046.7856
$ws.Range("I31").Value = 1152.5333
$ws.Range("J31").Value = 7386.3076
$ws.Range("K31").Value = 1152.5333
$ws.Range("L31").Value = 7386.3076
$ws.Range("M31").Value = -857.5333000000001
$ws.Range("N31").Value = -7976.3076
# Row 34 (hunk 23)
$ws.Range("H34").Value = 4046.7856
$ws.Range("I34").Value = 1152.5333
$ws.Range("J34").Value = 7386.3076
$ws.Range("K34").Value = 1152.5333
$ws.Range("L34").Value = 7386.3076
$ws.Range("M34").Value = -950.5333000000001
$ws.Range("N34").Value = -7790.3076
# Row 132 (hunk 24)
$ws.Range("H132").Value = 12823753
$ws.Range("I132").Value = 5363.6
$ws.Range("J132").Value = 20835246
$ws.Range("K132").Value = 16090.8
$ws.Range("L132").Value = 62505738
$ws.Range("M132").Value = -13560.8
$ws.Range("N132").Value = -62510798

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5 (hunk 25)
$ws.Range("H5").Value = 1090.6957
$ws.Range("I5").Value = 331.4
$ws.Range("J5").Value = 1674.7693
$ws.Range("K5").Value = 994.1999999999999
$ws.Range("L5").Value = 5024.3079
$ws.Range("M5").Value = -882.1999999999999
$ws.Range("N5").Value = -5248.3079
# Row 42 (hunk 26)
$ws.Range("H42").Value = 3166.6667
$ws.Range("I42").Value = 2000
$ws.Range("J42").Value = 3400
$ws.Range("K42").Value = 6000
$ws.Range("L42").Value = 10200
$ws.Range("M42").Value = -5466
$ws.Range("N42").Value = -11268
# Row 107 (hunk 27)
$ws.Range("H107").Value = 323.65714
$ws.Range("I107").Value = 300.48276
$ws.Range("J107").Value = 435.66666
$ws.Range("K107").Value = 901.44828
$ws.Range("L107").Value = 1306.99998
$ws.Range("M107").Value = 1018.55172
$ws.Range("N107").Value = -5146.999980000001
# Row 110 (hunk 28)
$ws.Range("H110").Value = 12502.875
$ws.Range("I110").Value = 3166.3333
$ws.Range("J110").Value = 13836.667
$ws.Range("K110").Value = 9498.999899999999
$ws.Range("L110").Value = 41510.001
$ws.Range("M110").Value = -5408.999899999999
$ws.Range("N110").Value = -49690.001
# Row 121 (hunk 29)
$ws.Range("H121").Value = 1131
$ws.Range("I121").Value = 454.42856
$ws.Range("J121").Value = 1279
$ws.Range("K121").Value = 1363.28568
$ws.Range("L121").Value = 3837
$ws.Range("M121").Value = -53.28567999999996
$ws.Range("N121").Value = -6457
# Row 135 (hunk 30)
$ws.Range("H135").Value = 1090.6957
$ws.Range("I135").Value = 331.4
$ws.Range("J135").Value = 1674.7693
$ws.Range("K135").Value = 2982.6
$ws.Range("L135").Value = 15072.9237
$ws.Range("M135").Value = -447.5999999999999
$ws.Range("N135").Value = -20142.9237
# Row 139 (hunk 31)
$ws.Range("H139").Value = 3308.525
$ws.Range("I139").Value = 1445.0952
$ws.Range("J139").Value = 5368.1055
$ws.Range("K139").Value = 4335.2856
$ws.Range("L139").Value = 16104.3165
$ws.Range("M139").Value = 804.7143999999998
$ws.Range("N139").Value = -26384.3165

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132 (hunk 32)
$ws.Range("H132").Value = 3814.2222
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 3041
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 9123
$ws.Range("M132").Value = -27470
$ws.Range("N132").Value = -14183

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132 (hunk 33)
$ws.Range("H132").Value = 2703.0193
$ws.Range("I132").Value = 2509.175
$ws.Range("J132").Value = 3349.1667
$ws.Range("K132").Value = 7527.525000000001
$ws.Range("L132").Value = 10047.5001
$ws.Range("M132").Value = -4997.525000000001
$ws.Range("N132").Value = -15107.5001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132 (hunk 34)
$ws.Range("H132").Value = 4506030.5
$ws.Range("I132").Value = 1751.2609
$ws.Range("J132").Value = 11905918
$ws.Range("K132").Value = 5253.7827
$ws.Range("L132").Value = 35717754
$ws.Range("M132").Value = -2723.7827
$ws.Range("N132").Value = -35722814
